# More QTI Picks #4s
# Replace the GE Aerospace "Buy" pick in row 5 with a new AMD "Short" pick,
# and add a new Pfizer "Short" pick in row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: was "Buy / GE Aerospace Inc. / GE / 193.68 / stop 187 / 10/21/24"
#     becomes "Short / AMD / AMD / 165.35 / stop 0 / 10/10/24"
$ws.Range("A5").Value = "Short"
$ws.Range("B5").Value = "AMD"
$ws.Range("C5").Value = "AMD"
$ws.Range("D5").Value = 165.35
$ws.Range("F5").Value = 0
# Day-bought column holds text like "10/10/24" - force Text format first so
# Excel doesn't auto-convert the string into a date serial number.
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "10/10/24"

# --- Row 7 (new row): "Short / Pfizer Inc. / PFE / 29.17 / stop 0 / 10/18/24"
$ws.Range("A7").Value = "Short"
$ws.Range("B7").Value = "Pfizer Inc."
$ws.Range("C7").Value = "PFE"
$ws.Range("D7").Value = 29.17
$ws.Range("F7").Value = 0
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "10/18/24"
# Match the normal "Day Bought" number format (matches the other rows) rather
# than leaving it on the plain Text format used only for the editing step above.
$ws.Range("H7").NumberFormat = "#,##0.00%"

# Move the active selection, mirroring where Excel left the cursor after
# entering the new data.
$ws.Range("G14").Select()
